# Scheduled market-data refresh: update cached price/profit figures
# across the per-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 237.61
$ws.Range("I15").Value = 237.61
$ws.Range("K15").Value = 712.83
$ws.Range("M15").Value = -543.83
$ws.Range("H32").Value = 871.0476
$ws.Range("J32").Value = 981.2308
$ws.Range("L32").Value = 981.2308
$ws.Range("N32").Value = -1633.2308
$ws.Range("H112").Value = 1301.9
$ws.Range("J112").Value = 1301.9
$ws.Range("L112").Value = 3905.7
$ws.Range("N112").Value = -6121.700000000001
$ws.Range("H129").Value = 1721.2307
$ws.Range("J129").Value = 1901.0217
$ws.Range("L129").Value = 5703.0651
$ws.Range("N129").Value = -15703.0651

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4293.7637
$ws.Range("I32").Value = 4688.02
$ws.Range("K32").Value = 4688.02
$ws.Range("M32").Value = -4401.02
$ws.Range("H35").Value = 22479.8
$ws.Range("I35").Value = 8400
$ws.Range("J35").Value = 31866.334
$ws.Range("K35").Value = 8400
$ws.Range("L35").Value = 31866.334
$ws.Range("M35").Value = -7994
$ws.Range("N35").Value = -32678.334
$ws.Range("H122").Value = 2804.375
$ws.Range("I122").Value = 1614.1666
$ws.Range("J122").Value = 6375
$ws.Range("K122").Value = 4842.4998
$ws.Range("L122").Value = 19125
$ws.Range("M122").Value = -2392.4998
$ws.Range("N122").Value = -24025
$ws.Range("H132").Value = 3347.182
$ws.Range("I132").Value = 2135.3076
$ws.Range("J132").Value = 4134.9
$ws.Range("K132").Value = 6405.9228
$ws.Range("L132").Value = 12404.7
$ws.Range("M132").Value = -3875.9228
$ws.Range("N132").Value = -17464.7

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1471.909
$ws.Range("I37").Value = 319.1
$ws.Range("J37").Value = 13000
$ws.Range("K37").Value = 319.1
$ws.Range("L37").Value = 13000
$ws.Range("M37").Value = -182.1
$ws.Range("N37").Value = -13274
$ws.Range("H95").Value = 32125
$ws.Range("J95").Value = 32125
$ws.Range("L95").Value = 32125
$ws.Range("N95").Value = -37617
$ws.Range("H103").Value = 35285.715
$ws.Range("J103").Value = 35285.715
$ws.Range("L103").Value = 35285.715
$ws.Range("N103").Value = -37629.715
$ws.Range("H134").Value = 3727.074
$ws.Range("I134").Value = 1643.1818
$ws.Range("J134").Value = 5159.75
$ws.Range("K134").Value = 4929.5454
$ws.Range("L134").Value = 15479.25
$ws.Range("M134").Value = -2394.5454
$ws.Range("N134").Value = -20549.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1975.6346
$ws.Range("I58").Value = 1099.4117
$ws.Range("J58").Value = 3630.7222
$ws.Range("K58").Value = 1099.4117
$ws.Range("L58").Value = 3630.7222
$ws.Range("M58").Value = -896.4117000000001
$ws.Range("N58").Value = -4036.7222
$ws.Range("H74").Value = 28329.3
$ws.Range("I74").Value = 5142.5
$ws.Range("J74").Value = 34126
$ws.Range("K74").Value = 5142.5
$ws.Range("L74").Value = 34126
$ws.Range("M74").Value = -4268.5
$ws.Range("N74").Value = -35874
$ws.Range("H77").Value = 28329.3
$ws.Range("I77").Value = 5142.5
$ws.Range("J77").Value = 34126
$ws.Range("K77").Value = 15427.5
$ws.Range("L77").Value = 102378
$ws.Range("M77").Value = -11059.5
$ws.Range("N77").Value = -111114
$ws.Range("H134").Value = 1408.6757
$ws.Range("I134").Value = 947.2273
$ws.Range("J134").Value = 2085.4666
$ws.Range("K134").Value = 2841.6819
$ws.Range("L134").Value = 6256.399800000001
$ws.Range("M134").Value = -306.6819
$ws.Range("N134").Value = -11326.3998
$ws.Range("H136").Value = 1975.6346
$ws.Range("I136").Value = 1099.4117
$ws.Range("J136").Value = 3630.7222
$ws.Range("K136").Value = 3298.2351
$ws.Range("L136").Value = 10892.1666
$ws.Range("M136").Value = -748.2351000000003
$ws.Range("N136").Value = -15992.1666
$ws.Range("H140").Value = 121083.75
$ws.Range("J140").Value = 121083.75
$ws.Range("L140").Value = 121083.75
$ws.Range("N140").Value = -131443.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 47687428
$ws.Range("J37").Value = 47687428
$ws.Range("L37").Value = 143062284
$ws.Range("N37").Value = -143062508
$ws.Range("H113").Value = 443
$ws.Range("I113").Value = 456.85715
$ws.Range("J113").Value = 428.62964
$ws.Range("K113").Value = 1370.57145
$ws.Range("L113").Value = 1285.88892
$ws.Range("M113").Value = 799.4285500000001
$ws.Range("N113").Value = -5625.888919999999
$ws.Range("H131").Value = 790.92786
$ws.Range("I131").Value = 456.66666
$ws.Range("J131").Value = 812.96704
$ws.Range("K131").Value = 1369.99998
$ws.Range("L131").Value = 2438.90112
$ws.Range("M131").Value = 3670.00002
$ws.Range("N131").Value = -12518.90112

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 1000000000
$ws.Range("I20").Value = 1000000000
$ws.Range("K20").Value = 1000000000
$ws.Range("M20").Value = -999999755
$ws.Range("H113").Value = 1349
$ws.Range("I113").Value = 1349
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1349
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 821
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 3972.8572
$ws.Range("I126").Value = 2990.5715
$ws.Range("J126").Value = 6428.5713
$ws.Range("K126").Value = 8971.7145
$ws.Range("L126").Value = 19285.7139
$ws.Range("M126").Value = -6501.7145
$ws.Range("N126").Value = -24225.7139

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2269.4119
$ws.Range("I22").Value = 1479.6666
$ws.Range("K22").Value = 1479.6666
$ws.Range("M22").Value = -1184.6666
$ws.Range("H27").Value = 2269.4119
$ws.Range("I27").Value = 1479.6666
$ws.Range("K27").Value = 1479.6666
$ws.Range("M27").Value = -1372.6666
$ws.Range("H46").Value = 1770
$ws.Range("I46").Value = 643.75
$ws.Range("J46").Value = 2300
$ws.Range("K46").Value = 643.75
$ws.Range("L46").Value = 2300
$ws.Range("M46").Value = -455.75
$ws.Range("N46").Value = -2676

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1185277.8
$ws.Range("I126").Value = 2580
$ws.Range("K126").Value = 7740
$ws.Range("M126").Value = -5270
$ws.Range("H132").Value = 7756707
$ws.Range("I132").Value = 6317.5
$ws.Range("K132").Value = 18952.5
$ws.Range("M132").Value = -16422.5
$ws.Range("H136").Value = 12105.6875
$ws.Range("I136").Value = 12676.777
$ws.Range("J136").Value = 11371.429
$ws.Range("K136").Value = 38030.331
$ws.Range("L136").Value = 34114.287
$ws.Range("M136").Value = -35480.331
$ws.Range("N136").Value = -39214.287

